$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

for ($r = 2; $r -le 259; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value()
    if ($v -eq $oldVal) {
        $cell.Value = $newVal
    }
}
